$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰" + [char]10 + `
    "✅ Dólar paralelo: 68" + [char]10 + `
    [char]10 + `
    "Binance" + [char]10 + `
    "✅ 1000 Bs = 1.66 = 6257.13 pesos" + [char]10 + `
    "✅ 6257.13 pesos = 1.66 = 877.56 Bs" + [char]10 + `
    [char]10 + `
    "Promedio competencia" + [char]10 + `
    "✅ Tasa pesos: 20" + [char]10 + `
    "✅ Tasa Bs: 20" + [char]10 + `
    "✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 603
$ws2.Range("O10").Value = 3773.05
$ws2.Range("N12").Value = 3778.99
$ws2.Range("O12").Value = 530
